# Updates the "Optical_Power" map workbook:
#  - E2 / E3 "Pendiente ADM" placeholders resolved to real OT numbers
#  - A brand new case (7755) inserted as the new first record in the
#    dated list, pushing the previously-row-28..86 records down one row
#    (rows 29..87)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Resolve the two pending OT numbers on rows 2 and 3 -----------------
# These columns are stored as text in the workbook (even though the
# content looks numeric), so the cell must be formatted as Text before
# the value is written, otherwise Excel would coerce it into a number.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "810804380"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "810804375"

# --- 2. Insert the new record as row 28, shifting the rest down ------------
$ws.Rows("28:28").Insert()

# Columns that must remain plain numbers (Attachments, Coordenada_X/Y)
$ws.Range("I28").Value = 1
$ws.Range("M28").Value = -58.453119
$ws.Range("N28").Value = -34.55489

# All remaining columns on the new row are text columns; every one of
# them (even the ones that look numeric, like "Caso", "Comuna" and "OT")
# must be stored as text, matching the rest of the sheet.
$textCols = @("A","B","C","D","E","F","G","H","J","K","L","O","P","Q","R")
foreach ($col in $textCols) {
    $ws.Range($col + "28").NumberFormat = "@"
}

$ws.Range("A28").Value = "7755"
$ws.Range("B28").Value = "10/24/2025"
$ws.Range("C28").Value = "Munich 1715"
$ws.Range("D28").Value = "13"
$ws.Range("E28").Value = "810447258"
$ws.Range("F28").Value = "Optical Power"
$ws.Range("G28").Value = "Pendiente"
$ws.Range("H28").Value = "Picada"
$ws.Range("J28").Value = "Cambio"
$ws.Range("K28").Value = "Sin equipos"
$ws.Range("L28").Value = "Pasante"
$ws.Range("O28").Value = "Saavedra"
$ws.Range("P28").Value = "Capital Norte"
$ws.Range("Q28").Value = "BLO-C"
$ws.Range("R28").Value = "Fuera de Poligono OVL"

# The row-insert operation materializes the handful of genuinely blank
# cells (previously rows 53/67/84, now 54/68/85 after the shift) as
# empty-string cells. Clear them back to true blanks to match the rest
# of the sheet.
$ws.Range("H54").Value = ""
$ws.Range("H68").Value = ""
$ws.Range("L85").Value = ""
